$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1926605504587156
$ws.Range("C2").Value = 0.536697247706422
$ws.Range("J2").Value = 0.04128440366972477
$ws.Range("P2").Value = 0.1192660550458716
$ws.Range("S2").Value = 0.1100917431192661
$ws.Range("B3").Value = 0.008130081300813009
$ws.Range("C3").Value = 0.03252032520325204
$ws.Range("J3").Value = 0.04878048780487805
$ws.Range("P3").Value = 0.7398373983739838
$ws.Range("S3").Value = 0.1707317073170732
$ws.Range("J4").Value = 0.103448275862069
$ws.Range("P4").Value = 0.5172413793103449
$ws.Range("S4").Value = 0.3793103448275862
$ws.Range("P5").Value = 0.75
$ws.Range("S5").Value = 0.25
$ws.Range("B6").Value = 0.03720930232558139
$ws.Range("D6").Value = 0.004651162790697674
$ws.Range("F6").Value = 0.04651162790697674
$ws.Range("J6").Value = 0.2325581395348837
$ws.Range("O6").Value = 0.04186046511627907
$ws.Range("Q6").Value = 0.1767441860465116
$ws.Range("R6").Value = 0.08372093023255814
$ws.Range("S6").Value = 0.3767441860465116
$ws.Range("B7").Value = 0.09090909090909091
$ws.Range("D7").Value = 0.03305785123966942
$ws.Range("E7").Value = 0.008264462809917356
$ws.Range("J7").Value = 0.06611570247933884
$ws.Range("O7").Value = 0.02479338842975207
$ws.Range("Q7").Value = 0.05785123966942149
$ws.Range("R7").Value = 0.04132231404958678
$ws.Range("S7").Value = 0.5867768595041323
$ws.Range("B8").Value = 0.09803921568627451
$ws.Range("D8").Value = 0.02450980392156863
$ws.Range("E8").Value = 0.002450980392156863
$ws.Range("F8").Value = 0.06617647058823529
$ws.Range("J8").Value = 0.0857843137254902
$ws.Range("O8").Value = 0.0196078431372549
$ws.Range("Q8").Value = 0.1372549019607843
$ws.Range("R8").Value = 0.08823529411764706
$ws.Range("S8").Value = 0.4779411764705883
$ws.Range("B9").Value = 0.08482142857142858
$ws.Range("D9").Value = 0.02232142857142857
$ws.Range("F9").Value = 0.08035714285714286
$ws.Range("J9").Value = 0.08482142857142858
$ws.Range("O9").Value = 0.008928571428571428
$ws.Range("Q9").Value = 0.1160714285714286
$ws.Range("R9").Value = 0.09821428571428571
$ws.Range("S9").Value = 0.5044642857142857
$ws.Range("B10").Value = 0.09387351778656126
$ws.Range("D10").Value = 0.01185770750988142
$ws.Range("E10").Value = 0.001976284584980237
$ws.Range("F10").Value = 0.08498023715415019
$ws.Range("J10").Value = 0.08498023715415019
$ws.Range("O10").Value = 0.01482213438735178
$ws.Range("Q10").Value = 0.1966403162055336
$ws.Range("R10").Value = 0.08596837944664032
$ws.Range("S10").Value = 0.424901185770751
$ws.Range("G11").Value = 0.1457286432160804
$ws.Range("J11").Value = 0.07537688442211055
$ws.Range("K11").Value = 0.1959798994974874
$ws.Range("L11").Value = 0.5678391959798995
$ws.Range("S11").Value = 0.01507537688442211
$ws.Range("G12").Value = 0.6386554621848739
$ws.Range("J12").Value = 0.2689075630252101
$ws.Range("L12").Value = 0.05042016806722689
$ws.Range("S12").Value = 0.04201680672268908
$ws.Range("G13").Value = 0.7407407407407407
$ws.Range("J13").Value = 0.1851851851851852
$ws.Range("S13").Value = 0.07407407407407407
$ws.Range("F15").Value = 0.0303030303030303
$ws.Range("H15").Value = 0.1515151515151515
$ws.Range("I15").Value = 0.07575757575757576
$ws.Range("J15").Value = 0.3535353535353535
$ws.Range("K15").Value = 0.09595959595959595
$ws.Range("M15").Value = 0.01515151515151515
$ws.Range("O15").Value = 0.04040404040404041
$ws.Range("S15").Value = 0.2373737373737374
$ws.Range("F16").Value = 0.01515151515151515
$ws.Range("H16").Value = 0.1666666666666667
$ws.Range("I16").Value = 0.1363636363636364
$ws.Range("J16").Value = 0.3712121212121212
$ws.Range("K16").Value = 0.1060606060606061
$ws.Range("M16").Value = 0.01515151515151515
$ws.Range("O16").Value = 0.07575757575757576
$ws.Range("S16").Value = 0.1136363636363636
$ws.Range("F17").Value = 0.01552795031055901
$ws.Range("H17").Value = 0.2142857142857143
$ws.Range("I17").Value = 0.1149068322981366
$ws.Range("J17").Value = 0.4254658385093168
$ws.Range("K17").Value = 0.05900621118012422
$ws.Range("M17").Value = 0.01552795031055901
$ws.Range("O17").Value = 0.04968944099378882
$ws.Range("S17").Value = 0.1055900621118012
$ws.Range("F18").Value = 0.03571428571428571
$ws.Range("H18").Value = 0.1488095238095238
$ws.Range("I18").Value = 0.1428571428571428
$ws.Range("J18").Value = 0.4404761904761905
$ws.Range("K18").Value = 0.07142857142857142
$ws.Range("M18").Value = 0.005952380952380952
$ws.Range("O18").Value = 0.06547619047619048
$ws.Range("S18").Value = 0.08928571428571429
$ws.Range("F19").Value = 0.01969178082191781
$ws.Range("H19").Value = 0.226027397260274
$ws.Range("I19").Value = 0.113013698630137
$ws.Range("J19").Value = 0.3690068493150685
$ws.Range("K19").Value = 0.08133561643835617
$ws.Range("M19").Value = 0.01626712328767123
$ws.Range("N19").Value = 0.0008561643835616438
$ws.Range("O19").Value = 0.07534246575342465
